# Auto-generated PowerShell Excel COM-interop edit script
# Applies the cell-value updates from the commit diff to Results_MackeyGlass sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.01987698815834519
$ws.Range("D3").Value = 0.07747178101776631
$ws.Range("E3").Value = 0.009604327749957774
$ws.Range("G3").Value = 'max\_depth: 16, max\_features: 4 \\'
$ws.Range("C4").Value = 0.01002701242761721
$ws.Range("D4").Value = 0.03908089620351484
$ws.Range("E4").Value = 0.006632192080453148
$ws.Range("G4").Value = 'max\_depth: 16, max\_features: 8, n\_estimators: 200 \\'
$ws.Range("E5").Value = 0.06243720571138975
$ws.Range("C6").Value = 0.0998254113003442
$ws.Range("D6").Value = 0.3890756659238543
$ws.Range("E6").Value = 0.08762964479374162
$ws.Range("B7").Value = 'GBM \cite{friedman2001greedy}'
$ws.Range("C7").Value = 0.009690697544606516
$ws.Range("D7").Value = 0.03777008831038435
$ws.Range("E7").Value = 0.005494791757575131
$ws.Range("G7").Value = 'learning\_rate: 0.05, max\_depth: 20, max\_features: 2, n\_estimators: 200 \\'
$ws.Range("C10").Value = 0.04617939703653087
$ws.Range("D10").Value = 0.1799870335609462
$ws.Range("E10").Value = 0.03244929935914024
$ws.Range("C11").Value = 0.02283476574078957
$ws.Range("D11").Value = 0.08899990063734688
$ws.Range("E11").Value = 0.01819140053886059
$ws.Range("C12").Value = 0.01709613338723205
$ws.Range("D12").Value = 0.06663322891149877
$ws.Range("E12").Value = 0.01262624834438492
$ws.Range("C13").Value = 0.0131910578530774
$ws.Range("D13").Value = 0.05141295739803758
$ws.Range("E13").Value = 0.009526480523523931
$ws.Range("C14").Value = 0.0192742790299282
$ws.Range("D14").Value = 0.07512268520696433
$ws.Range("E14").Value = 0.01487175495059627
$ws.Range("C15").Value = 0.0209119019797498
$ws.Range("D15").Value = 0.08150542113997246
$ws.Range("E15").Value = 0.01700731925401119
$ws.Range("C16").Value = 0.09971232019375161
$ws.Range("D16").Value = 0.3886348863965334
$ws.Range("E16").Value = 0.0874389428291678
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 'omega: 50, r: 0.1 \\'
$ws.Range("B17").Value = 'Simpl\_eTS \cite{angelov2005simpl_ets}'
$ws.Range("C17").Value = 0.09971232019375161
$ws.Range("D17").Value = 0.3886348863965334
$ws.Range("E17").Value = 0.0874389428291678
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 'omega: 50, r: 0.1 \\'
$ws.Range("C18").Value = 0.08970967132845216
$ws.Range("D18").Value = 0.3496489486721231
$ws.Range("E18").Value = 0.07612164068576284
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 'mu: 0.5, omega: 1000 \\'
$ws.Range("C19").Value = 0.09703642017414169
$ws.Range("D19").Value = 0.3782054018743716
$ws.Range("E19").Value = 0.08355445192516217
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 'alpha: 0.001, beta: 0.2, lambda1: 0.001, omega: 10000, sigma: 0.25 \\'
$ws.Range("C20").Value = 0.03886467419850697
$ws.Range("D20").Value = 0.1514774524615017
$ws.Range("E20").Value = 0.03471813511043065
$ws.Range("F20").Value = 74
$ws.Range("G20").Value = 'alpha: 0.001, lambda1: 0.5, omega: 10000, sigma: 0.003, w: 50 \\'
$ws.Range("C21").Value = 0.09959557561220603
$ws.Range("D21").Value = 0.3881798672264027
$ws.Range("E21").Value = 0.08732713427523939
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 'alpha: 0.001, beta: 0.01, e\_utility: 0.03, lambda1: 0.25, omega: 100, pi: 0.3, sigma: 0.1 \\'
$ws.Range("C22").Value = 0.005907353391721347
$ws.Range("D22").Value = 0.0230242723249726
$ws.Range("E22").Value = 0.00384480535989134
$ws.Range("F22").Value = 26
$ws.Range("G22").Value = 'alpha: 0.1, beta: 0.1, e\_utility: 0.03, lambda1: 1e-07, sigma: 0.5 \\'
$ws.Range("C23").Value = 0.1009652268874646
$ws.Range("D23").Value = 0.3935181671148084
$ws.Range("E23").Value = 0.08516756904634615
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 'fuzzy\_operator: prod, rules: 16 \\'
$ws.Range("C24").Value = 0.09954526595210997
$ws.Range("D24").Value = 0.3879837822391304
$ws.Range("E24").Value = 0.08515069302675088
$ws.Range("C25").Value = 0.09096224060842921
$ws.Range("D25").Value = 0.3545309142996582
$ws.Range("E25").Value = 0.07995993584213755
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 'adaptive\_filter: wRLS, fuzzy\_operator: prod, rules: 17 \\'
$ws.Range("C26").Value = 0.1022941964254241
$ws.Range("D26").Value = 0.3986979074358216
$ws.Range("E26").Value = 0.08748953164459831
$ws.Range("F26").Value = 13
$ws.Range("G26").Value = 'error\_metric: MAE, fuzzy\_operator: minmax, num\_generations: 10, num\_parents\_mating: 5, parallel\_processing: 10, rules: 13, sol\_per\_pop: 10 \\'
$ws.Range("C27").Value = 0.1068393136437876
$ws.Range("D27").Value = 0.4164127806870441
$ws.Range("E27").Value = 0.0866251704022159
$ws.Range("G27").Value = 'adaptive\_filter: RLS, error\_metric: CPPM, fuzzy\_operator: prod, lambda1: 0.97, num\_generations: 5, num\_parents\_mating: 5, parallel\_processing: 10, rules: 1, sol\_per\_pop: 5 \\'
$ws.Range("C28").Value = 0.09096224060842921
$ws.Range("D28").Value = 0.3545309142996582
$ws.Range("E28").Value = 0.07995993584213755
$ws.Range("F28").Value = 17
$ws.Range("G28").Value = 'adaptive\_filter: wRLS, error\_metric: MAE, fuzzy\_operator: prod, num\_generations: 5, num\_parents\_mating: 5, parallel\_processing: 10, rules: 17, sol\_per\_pop: 5 \\'
$ws.Range("C29").Value = 0.1048791579970679
$ws.Range("D29").Value = 0.408772953776967
$ws.Range("E29").Value = 0.08916767565294725
$ws.Range("C30").Value = 0.09741149189363492
$ws.Range("D30").Value = 0.3796672669158486
$ws.Range("E30").Value = 0.08524810957110011
$ws.Range("C31").Value = 0.01482123194899861
$ws.Range("D31").Value = 0.05776666096590034
$ws.Range("E31").Value = 0.0122677759938756
